$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 18 ("Quorum-based" slide): merge split runs that have identical
# formatting back into single runs. Re-assigning the same text to a
# Characters() sub-range spanning the two original runs collapses them into
# one run (same behaviour as PowerPoint's "remove extra run" autocleanup).
# ---------------------------------------------------------------------------
$s18 = $p.Slides.Item(18)
$sh18 = $s18.Shapes.Item(2)
$tr18 = $sh18.TextFrame.TextRange

# "系统的" + "问题空间" -> "系统的问题空间"
$rng = $tr18.Characters(53, 7)
$rng.Text = $rng.Text

# "系统及调研相关" + "工具" -> "系统及调研相关工具"
$rng = $tr18.Characters(82, 9)
$rng.Text = $rng.Text

# ---------------------------------------------------------------------------
# Slide 22 (Cassandra physical distribution): "个数据" + "中心" -> "个数据中心"
# ---------------------------------------------------------------------------
$s22 = $p.Slides.Item(22)
$sh22 = $s22.Shapes.Item(2)
$tr22 = $sh22.TextFrame.TextRange
$rng = $tr22.Characters(34, 5)
$rng.Text = $rng.Text

# ---------------------------------------------------------------------------
# Slide 23 (Cassandra running environment): shrink the oversized textbox and
# merge the split runs.
# ---------------------------------------------------------------------------
$s23 = $p.Slides.Item(23)
$sh23 = $s23.Shapes.Item(2)
$tr23 = $sh23.TextFrame.TextRange

# "远程" + "批量" + "启动和停止" -> "远程批量启动和停止"
$rng = $tr23.Characters(283, 9)
$rng.Text = $rng.Text

# "远程运行日志" + "汇总" -> "远程运行日志汇总"
$rng = $tr23.Characters(302, 8)
$rng.Text = $rng.Text

# fix oversized textbox: cy 7478970 -> 5024218 EMU (914400 EMU/in, 12700 EMU/pt)
$sh23.Height = 5024218 / 12700

# ---------------------------------------------------------------------------
# Slide 28 (inversion idea): "权值高，不同" + "进程" -> "权值高，不同进程"
# ---------------------------------------------------------------------------
$s28 = $p.Slides.Item(28)
$sh28 = $s28.Shapes.Item(2)
$tr28 = $sh28.TextFrame.TextRange
$rng = $tr28.Characters(135, 8)
$rng.Text = $rng.Text

# ---------------------------------------------------------------------------
# Slide 29 (inversion example): fix oversized textbox cy 3471720 -> 2304256
# ---------------------------------------------------------------------------
$s29 = $p.Slides.Item(29)
$sh29 = $s29.Shapes.Item(5)
$sh29.Height = 2304256 / 12700
